$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new test-results row for krb5.cap, matching the look of the other
# "tested" rows (copy row 17's formatting/row-height, then overwrite values)
$ws.Range("A17:D17").Copy($ws.Range("A25:D25"))
$ws.Rows.Item(25).RowHeight = $ws.Rows.Item(17).RowHeight

$ws.Range("A25").Value = "krb5.cap"
$ws.Range("B25").Value = "msg_type, cname, sname, err_code, weak_encryption, ticket_encryption"
$ws.Range("C25").Value = "112-7"
$ws.Range("D25").Value = "yes"

# The view had scrolled up one row by the time this row was added
$ws.Application.ActiveWindow.ScrollRow = 18
